# This script updates the NATMI LR-pairs sheet (Wnt5a-Fzd7) with refreshed TPM-derived
# values for the "Receptor average/total expression value" columns (M, N) and the
# downstream columns that are recomputed from them:
#   O = Receptor derived specificity of average expression value
#       = M / (sum of M over the 5 target clusters within the same Sending-cluster group)
#   P = Receptor derived specificity of total expression value
#       = N / (sum of N over the 5 target clusters within the same Sending-cluster group)
#   Q = Edge average expression weight        = (Ligand avg expr, col G) * M
#   R = Edge total expression weight          = (Ligand total expr, col H) * N
#   S = Edge average expression derived specificity = Q / (sum of Q over all 10 data rows)
#   T = Edge total expression derived specificity   = R / (sum of R over all 10 data rows)
#
# Only the rows/columns that actually change per the updated TPM values are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0.7339303333333334
$ws.Range("N2").Value = 2.201791
$ws.Range("O2").Value = 0.03574007706012852
$ws.Range("P2").Value = 0.03574007706012852
$ws.Range("Q2").Value = 0.1173557049434445
$ws.Range("R2").Value = 1.056201344491
$ws.Range("S2").Value = 0.0009463157329718295
$ws.Range("T2").Value = 0.0009463157329718295
$ws.Range("O3").Value = 0.3842514532634088
$ws.Range("P3").Value = 0.3842514532634088
$ws.Range("S3").Value = 0.01017410217187556
$ws.Range("T3").Value = 0.01017410217187556
$ws.Range("M4").Value = 4.974008666666667
$ws.Range("N4").Value = 14.922026
$ws.Range("O4").Value = 0.2422184299659874
$ws.Range("P4").Value = 0.2422184299659874
$ws.Range("Q4").Value = 0.795345643802889
$ws.Range("R4").Value = 7.158110794226
$ws.Range("S4").Value = 0.006413391630547449
$ws.Range("T4").Value = 0.006413391630547449
$ws.Range("M5").Value = 2.087648
$ws.Range("N5").Value = 6.262943999999999
$ws.Range("O5").Value = 0.1016618294757629
$ws.Range("P5").Value = 0.1016618294757629
$ws.Range("Q5").Value = 0.3338156110826666
$ws.Range("R5").Value = 3.004340499744
$ws.Range("S5").Value = 0.002691773398075258
$ws.Range("T5").Value = 0.002691773398075258
$ws.Range("M6").Value = 4.848944666666667
$ws.Range("N6").Value = 14.546834
$ws.Range("O6").Value = 0.2361282102347124
$ws.Range("P6").Value = 0.2361282102347124
$ws.Range("Q6").Value = 0.775347868514889
$ws.Range("R6").Value = 6.978130816634
$ws.Range("S6").Value = 0.006252136501207213
$ws.Range("T6").Value = 0.006252136501207213
$ws.Range("M7").Value = 0.7339303333333334
$ws.Range("N7").Value = 2.201791
$ws.Range("O7").Value = 0.03574007706012852
$ws.Range("P7").Value = 0.03574007706012852
$ws.Range("Q7").Value = 4.314887987077333
$ws.Range("R7").Value = 38.833991883696
$ws.Range("S7").Value = 0.03479376132715669
$ws.Range("T7").Value = 0.03479376132715669
$ws.Range("O8").Value = 0.3842514532634088
$ws.Range("P8").Value = 0.3842514532634088
$ws.Range("S8").Value = 0.3740773510915332
$ws.Range("T8").Value = 0.3740773510915332
$ws.Range("M9").Value = 4.974008666666667
$ws.Range("N9").Value = 14.922026
$ws.Range("O9").Value = 0.2422184299659874
$ws.Range("P9").Value = 0.2422184299659874
$ws.Range("Q9").Value = 29.24295300065067
$ws.Range("R9").Value = 263.186577005856
$ws.Range("S9").Value = 0.2358050383354399
$ws.Range("T9").Value = 0.2358050383354399
$ws.Range("M10").Value = 2.087648
$ws.Range("N10").Value = 6.262943999999999
$ws.Range("O10").Value = 0.1016618294757629
$ws.Range("P10").Value = 0.1016618294757629
$ws.Range("Q10").Value = 12.273599914496
$ws.Range("R10").Value = 110.462399230464
$ws.Range("S10").Value = 0.09897005607768766
$ws.Range("T10").Value = 0.09897005607768766
$ws.Range("M11").Value = 4.848944666666667
$ws.Range("N11").Value = 14.546834
$ws.Range("O11").Value = 0.2361282102347124
$ws.Range("P11").Value = 0.2361282102347124
$ws.Range("Q11").Value = 28.50768273492267
$ws.Range("R11").Value = 256.569144614304
$ws.Range("S11").Value = 0.2298760737335052
$ws.Range("T11").Value = 0.2298760737335052
